$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Price (D) and Volume(1h) (E) values per the Jan 28 2023 symbol-list refresh.
# Leading apostrophe keeps these numeric-looking strings stored as text (matching
# the source data, which is textual, not numeric).
$updates = @(
    @{ Row = 2; D = "307.34"; E = "0.78%" }
    @{ Row = 3; D = "38.84"; E = "8.08%" }
    @{ Row = 4; D = "5.089"; E = "0.91%" }
    @{ Row = 5; D = "0.08062"; E = "0.34%" }
    @{ Row = 6; D = "1.921"; E = "2.75%" }
    @{ Row = 7; D = "4.186"; E = "1.62%" }
    @{ Row = 8; D = "7.937"; E = "1.93%" }
    @{ Row = 9; D = "0.9304"; E = "0.45%" }
    @{ Row = 10; D = "0.1446"; E = "10.55%" }
    @{ Row = 11; D = "0.1938"; E = "2.36%" }
    @{ Row = 12; D = "0.08955"; E = "-1.18%" }
    @{ Row = 13; D = "0.03495"; E = "1.66%" }
    @{ Row = 14; D = "0.09781"; E = "-0.79%" }
    @{ Row = 15; D = "0.001401"; E = "-0.20%" }
    @{ Row = 16; D = "0.005989"; E = "-2.93%" }
    @{ Row = 17; D = "3.750"; E = "-2.23%" }
    @{ Row = 18; D = "3.464"; E = "2.15%" }
    @{ Row = 19; D = "0.3445"; E = "0.91%" }
    @{ Row = 20; D = "0.1320"; E = "-0.94%" }
    @{ Row = 21; D = "4.776"; E = "-0.69%" }
    @{ Row = 22; D = $null; E = "4.90%" }
    @{ Row = 23; D = "0.04375"; E = "0.19%" }
    @{ Row = 24; D = "0.001234"; E = "0.40%" }
    @{ Row = 25; D = "0.004278"; E = "-0.13%" }
    @{ Row = 26; D = "0.0001299"; E = "-0.03%" }
    @{ Row = 39; D = "0.02071"; E = "3.80%" }
    @{ Row = 40; D = "0.05068"; E = "-0.96%" }
    @{ Row = 41; D = "0.007386"; E = "-1.52%" }
    @{ Row = 42; D = "0.010000"; E = "-0.57%" }
    @{ Row = 43; D = $null; E = "-0.28%" }
    @{ Row = 44; D = "0.002138"; E = "-1.41%" }
    @{ Row = 45; D = "0.009071"; E = "-5.58%" }
    @{ Row = 46; D = "0.00006172"; E = "-0.51%" }
    @{ Row = 47; D = $null; E = "0.10%" }
    @{ Row = 48; D = "0.002797"; E = $null }
    @{ Row = 49; D = "0.001598"; E = "28.05%" }
    @{ Row = 50; D = "0.00002099"; E = "0.10%" }
    @{ Row = 51; D = "0.0001999"; E = "0.10%" }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $ws.Cells.Item($u.Row, 4).Value = "'" + $u.D
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($u.Row, 5).Value = "'" + $u.E
    }
}
